# repull data, push all data, mean calculation
# Update column F (dSF) values to match repulled data for select rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -8
$ws.Range("F5").Value  = -8
$ws.Range("F7").Value  = -3
$ws.Range("F9").Value  = -8
$ws.Range("F10").Value = -8
$ws.Range("F11").Value = -8
$ws.Range("F13").Value = -2
$ws.Range("F14").Value = -3
$ws.Range("F16").Value = -6
$ws.Range("F17").Value = 12
$ws.Range("F19").Value = 3
